$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price/Volume columns keep their text representation (no numeric coercion)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = '29.373.60'
$ws.Range("E2").Value = '  -0.78%  '
$ws.Range("D3").Value = '1.846.91'
$ws.Range("E3").Value = '  -0.43%  '
$ws.Range("D4").Value = '0.9987'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").Value = '241.27'
$ws.Range("E5").Value = '  -1.01%  '
$ws.Range("D6").Value = '0.6330'
$ws.Range("E6").Value = '  -1.22%  '
$ws.Range("E7").Value = '  +0.02%  '
$ws.Range("D8").Value = '0.07565'
$ws.Range("E8").Value = '  -0.34%  '
$ws.Range("D9").Value = '0.2961'
$ws.Range("E9").Value = '  -1.43%  '
$ws.Range("D10").Value = '24.86'
$ws.Range("E10").Value = '  +1.77%  '
$ws.Range("D11").Value = '0.07739'
$ws.Range("E11").Value = '  +0.73%  '
$ws.Range("D12").Value = '5.001'
$ws.Range("E12").Value = '  -1.11%  '
$ws.Range("D13").Value = '0.6825'
$ws.Range("E13").Value = '  -1.17%  '
$ws.Range("D14").Value = '83.01'
$ws.Range("E14").Value = '  -1.49%  '
$ws.Range("D15").Value = '0.000009964'
$ws.Range("E15").Value = '  +2.87%  '
$ws.Range("D16").Value = '6.138'
$ws.Range("E16").Value = '  -2.85%  '
$ws.Range("D17").Value = '29.417.21'
$ws.Range("E17").Value = '  -0.76%  '
$ws.Range("D18").Value = '230.41'
$ws.Range("E18").Value = '  -3.75%  '
$ws.Range("E19").Value = '  -1.36%  '
$ws.Range("D20").Value = '0.9999'
$ws.Range("E20").Value = '  +0.00%  '
$ws.Range("D21").Value = '7.549'
$ws.Range("E21").Value = '  -1.26%  '
$ws.Range("E23").Value = '  +232.11%  '
$ws.Range("D24").Value = '16.64'
$ws.Range("E24").Value = '  +171.06%  '
$ws.Range("D25").Value = '156.46'
$ws.Range("E25").Value = '  -0.41%  '
$ws.Range("D26").Value = '0.1398'
$ws.Range("E26").Value = '  -0.60%  '
$ws.Range("D27").Value = '8.386'
$ws.Range("E27").Value = '  -1.55%  '
$ws.Range("D28").Value = '17.68'
$ws.Range("E28").Value = '  -0.65%  '
$ws.Range("D29").Value = '2.714'
$ws.Range("E29").Value = '  +171.82%  '
$ws.Range("D30").Value = '1.469'
$ws.Range("E30").Value = '  -1.39%  '
$ws.Range("D31").Value = '0.05704'
$ws.Range("E31").Value = '  -3.13%  '
$ws.Range("E32").Value = '  -2.45%  '
$ws.Range("D33").Value = '4.128'
$ws.Range("E33").Value = '  -0.52%  '
$ws.Range("D34").Value = '4.020'
$ws.Range("E34").Value = '  -1.55%  '
$ws.Range("D35").Value = '1.845'
$ws.Range("E35").Value = '  -3.40%  '
$ws.Range("D36").Value = '1.154'
$ws.Range("E36").Value = '  -2.74%  '
$ws.Range("D37").Value = '0.7164'
$ws.Range("E37").Value = '  -1.19%  '
$ws.Range("E38").Value = '  +0.02%  '
$ws.Range("D39").Value = '1.244.03'
$ws.Range("E39").Value = '  +2.15%  '
$ws.Range("D40").Value = '2.805'
$ws.Range("E40").Value = '  -0.07%  '
$ws.Range("D41").Value = '0.01807'
$ws.Range("E41").Value = '  +1.53%  '
$ws.Range("E42").Value = '  +264.78%  '
$ws.Range("D43").Value = '0.9009'
$ws.Range("E43").Value = '  -1.74%  '
$ws.Range("D44").Value = '0.9997'
$ws.Range("E44").Value = '  -0.02%  '
$ws.Range("D45").Value = '101.90'
$ws.Range("E45").Value = '  -0.05%  '
$ws.Range("D46").Value = '66.18'
$ws.Range("E46").Value = '  -2.23%  '
$ws.Range("B47").Value = 'Aptos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D47").Value = '7.051'
$ws.Range("E47").Value = '  -6.08%  '
$ws.Range("B48").Value = 'EnergySwap'
$ws.Range("C48").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D48").Value = '9.147'
$ws.Range("E48").Value = '  -0.94%  '
$ws.Range("B49").Value = 'TheSandbox'
$ws.Range("C49").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D49").Value = '0.4021'
$ws.Range("E49").Value = '  -1.20%  '
$ws.Range("B50").Value = 'RenderToken'
$ws.Range("C50").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D50").Value = '1.704'
$ws.Range("E50").Value = '  +1.03%  '
$ws.Range("D51").Value = '0.1125'
$ws.Range("E51").Value = '  -0.66%  '
